$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date formatting from D1 into D5 so the new row reuses the
# existing "short date" cell style instead of creating a brand new one.
$ws.Cells.Item(1, 4).Copy()
$ws.Cells.Item(5, 4).PasteSpecial(-4122)  # xlPasteFormats

# Populate the new 5th data row: Anna Lushnikova, same date as row 1,
# line "1", station "Belorusskaya" (same as row 1).
$ws.Cells.Item(5, 1).Value = 5
$ws.Cells.Item(5, 2).Value = "Anna"
$ws.Cells.Item(5, 3).Value = "Lushnikova"
$ws.Cells.Item(5, 4).Value = 45292
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = "Belorusskaya"

# Move the active selection to G6, matching the saved view state.
[void]$ws.Range("G6").Select()
